# Sync the event-listing sheets ("展览" and "全部类型") to the next scrape
# snapshot: the oldest event (2024-03-23 "Look Look动漫嘉年华") has dropped off
# the list, every remaining event's row shifts up by one, a handful of
# "想去人数" (interest count) values tick up/down slightly, and the table
# shrinks from 29 data rows to 28 (dimension A1:I29 -> A1:I28). The leading
# index column (A) is left untouched -- only columns B..I move.

$rowsData = @(
    ,@("2024-03-30", "苏州·Anime Space动漫游戏展", "金芳路与新发路交叉口东南120米 万龙大厦", "2024.03.30 10:00-03.30 17:00", 614, 55, "https://show.bilibili.com/platform/detail.html?id=82815", "//i0.hdslb.com/bfs/openplatform/202403/XPHUZMJa1710327274324.jpeg")
    ,@("2024-03-30", "苏州·奇幻世界5.3动漫游戏展", "龙河路1288号 乐动力苏州湾体育中心", "2024.03.30 10:00-03.31 17:00", 2181, 55, "https://show.bilibili.com/platform/detail.html?id=82002", "//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg")
    ,@("2024-03-31", "张家港·META萌圆饿了", "大新镇人民路18号 新香苑宴会厅", "2024.03.31 10:00-03.31 17:00", 71, 30, "https://show.bilibili.com/platform/detail.html?id=82407", "//i0.hdslb.com/bfs/openplatform/202403/pxZkSPcL1709707210340.jpeg")
    ,@("2024-04-04", "【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会", "金山南路288号木渎影视城F2 苏州广电国际会展中心", "2024.04.04 10:00-04.05 17:00", 12775, 60, "https://show.bilibili.com/platform/detail.html?id=81827", "//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg")
    ,@("2024-04-04", "常熟·SL动漫展", "报慈北路218号 四季花园酒店(报慈北路店)", "2024.04.04 10:00-04.04 17:00", 63, 50, "https://show.bilibili.com/platform/detail.html?id=82250", "//i0.hdslb.com/bfs/openplatform/202403/rcvwgj7N1709518723294.jpeg")
    ,@("2024-04-06", "苏州·第一届寒假动漫展宅舞比赛-CF01", "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店", "2024.04.06 10:00-04.06 16:00", 112, 49, "https://show.bilibili.com/platform/detail.html?id=80528", "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg")
    ,@("2024-04-13", "苏州·X-party 国漫游戏嘉年华03", "秋枫街与开平路交叉口西南角 爱琴海购物中心", "2024.04.13 10:00-04.14 17:00", 506, 48, "https://show.bilibili.com/platform/detail.html?id=82042", "//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg")
    ,@("2024-04-13", "苏州·绘时国乙1.0-秩序之外", "兴中路与鲈乡北路交汇处 香漫商业广场", "2024.04.13 13:00-04.13 20:00", 467, 88, "https://show.bilibili.com/platform/detail.html?id=80789", "//i2.hdslb.com/bfs/openplatform/202403/DI2ackIO1710137864319.jpeg")
    ,@("2024-04-20", "苏州·首届Redamancy动漫游戏嘉年华", "清禾路886号 尹山湖大剧院", "2024.04.20 10:00-04.20 17:00", 1159, 60, "https://show.bilibili.com/platform/detail.html?id=81879", "//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg")
    ,@("2024-04-21", "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0", "清禾路888号2号楼3楼 格莱美婚礼宴会中心", "2024.04.21 10:00-04.21 21:00", 960, 69.90000000000001, "https://show.bilibili.com/platform/detail.html?id=78666", "//i0.hdslb.com/bfs/openplatform/202403/PlZCFPVs1710502485559.jpeg")
    ,@("2024-05-01", "昆山·第十二届理想乡动漫游戏展", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.01 10:00-05.03 17:00", 13687, 75, "https://show.bilibili.com/platform/detail.html?id=77196", "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png")
    ,@("2024-05-01", "苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区", "金山南路288号 广电国际会展中心", "2024.05.01 10:00-05.02 17:00", 14105, 65, "https://show.bilibili.com/platform/detail.html?id=79789", "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg")
    ,@("2024-05-02", "昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.02 14:00-05.02 16:00", 40, 1, "https://show.bilibili.com/platform/detail.html?id=81116", "//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg")
    ,@("2024-05-02", "昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.02 14:00-05.02 16:00", 168, 1, "https://show.bilibili.com/platform/detail.html?id=81100", "//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg")
    ,@("2024-05-02", "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.02 14:00-05.02 16:00", 21, 1, "https://show.bilibili.com/platform/detail.html?id=81119", "//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg")
    ,@("2024-05-02", "昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.02 14:00-05.02 16:00", 42, 1, "https://show.bilibili.com/platform/detail.html?id=81118", "//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg")
    ,@("2024-05-02", "苏州·动漫游戏嘉年华", "东太湖大道12000号 苏州湾大剧院", "2024.05.02 10:00-05.02 17:00", 22, 58, "https://show.bilibili.com/platform/detail.html?id=82824", "//i1.hdslb.com/bfs/openplatform/202403/HzWBEJeE1710324788092.jpeg")
    ,@("2024-05-02", "苏州·动漫游戏嘉年华cv刘圣博见面会（取消）", "东太湖大道12000号 苏州湾大剧院", "2024.05.02 10:00-05.02 17:00", 1, "不可售", "https://show.bilibili.com/platform/detail.html?id=83038", "//i0.hdslb.com/bfs/openplatform/202403/D94B39u21710901393375.jpeg")
    ,@("2024-05-02", "苏州·动漫游戏嘉年华cv张文钰见面会（取消）", "东太湖大道12000号 苏州湾大剧院", "2024.05.02 10:00-05.02 17:00", 1, "不可售", "https://show.bilibili.com/platform/detail.html?id=83037", "//i1.hdslb.com/bfs/openplatform/202403/eP5thEuS1710901472569.jpeg")
    ,@("2024-05-02", "苏州·动漫游戏嘉年华cv沐霏见面会（取消）", "东太湖大道12000号 苏州湾大剧院", "2024.05.02 10:00-05.02 17:00", 4, "不可售", "https://show.bilibili.com/platform/detail.html?id=82891", "//i1.hdslb.com/bfs/openplatform/202403/8VORpvQc1710900704258.jpeg")
    ,@("2024-05-03", "常熟·CDW·动漫展03", "常熟国际展览中心 国际展览中心", "2024.05.03 09:00-05.04 17:30", 1067, 60, "https://show.bilibili.com/platform/detail.html?id=82489", "//i0.hdslb.com/bfs/openplatform/202403/XK411blC1709794808211.jpeg")
    ,@("2024-05-03", "昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.03 14:00-05.03 16:00", 109, 1, "https://show.bilibili.com/platform/detail.html?id=81120", "//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg")
    ,@("2024-05-03", "昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会", "花桥经济开发区绿地大道1598号 花桥国际博览中心", "2024.05.03 14:00-05.03 16:00", 56, 1, "https://show.bilibili.com/platform/detail.html?id=81114", "//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg")
    ,@("2024-05-04", "【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华阿杰内场", "苏州大道东688号 苏州国际博览中心", "2024.05.04 09:00-05.04 17:00", 933, "已售罄", "https://show.bilibili.com/platform/detail.html?id=82940", "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg")
    ,@("2024-05-04", "苏州·OCG国潮动漫游戏嘉年华", "苏州大道东688号 苏州国际博览中心", "2024.05.04 09:00-05.05 17:00", 5215, 65, "https://show.bilibili.com/platform/detail.html?id=82779", "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg")
    ,@("2024-05-18", "苏州·YoungComic动漫嘉年华", "清禾路886号 尹山湖大剧院", "2024.05.18 10:00-05.18 17:00", 8, 60, "https://show.bilibili.com/platform/detail.html?id=83142", "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg")
    ,@("2024-06-08", "【会员购严选】苏州·Come in joy动漫国潮文化节", "金山南路288号 广电国际会展中心", "2024.06.08 10:00-06.09 17:00", 271, 60, "https://show.bilibili.com/platform/detail.html?id=82233", "//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg")
)

function Update-EventSheet {
    param($ws)

    # Column B holds plain-text dates like "2024-03-30". Left alone, Excel's
    # COM layer auto-recognizes that shape and silently coerces it to a real
    # date serial, which is not what the source feed stores. Pre-formatting
    # the column as text keeps the literal string.
    $ws.Range("B2:B28").NumberFormat = "@"

    for ($idx = 0; $idx -lt $rowsData.Count; $idx++) {
        $r = $idx + 2   # data starts on row 2 (row 1 is the header)
        $vals = $rowsData[$idx]
        $ws.Cells.Item($r, 2).Value = $vals[0]   # B 开始时间
        $ws.Cells.Item($r, 3).Value = $vals[1]   # C 名称
        $ws.Cells.Item($r, 4).Value = $vals[2]   # D 地点
        $ws.Cells.Item($r, 5).Value = $vals[3]   # E 具体时间范围
        $ws.Cells.Item($r, 6).Value = $vals[4]   # F 想去人数
        $ws.Cells.Item($r, 7).Value = $vals[5]   # G 最低票价
        $ws.Cells.Item($r, 8).Value = $vals[6]   # H Link
        $ws.Cells.Item($r, 9).Value = $vals[7]   # I Cover
    }

    # The old row 29 (previously the last data row) is no longer part of the
    # feed; drop it so the sheet shrinks back down to A1:I28.
    $ws.Rows.Item(29).Delete()
}

$wb = $excel.ActiveWorkbook
Update-EventSheet $wb.Worksheets.Item("展览")
Update-EventSheet $wb.Worksheets.Item("全部类型")
